# Apply the "Polish README and update dependencies" content edit to the
# sample AI deck: delete the References slide, rewrite the remaining four
# slides with AI-themed copy, and add speaker notes to the first three.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Helper: set a single paragraph's text robustly. Some placeholder runs
# (notably ones that already start with the literal word "Item") trip an
# autocorrect-style run split when overwritten directly (e.g. "Involves"
# becomes separate "I" / "nvolves" runs). Writing a neutral placeholder
# first and then the real text sidesteps that quirk.
# ---------------------------------------------------------------------
function Set-ParaText($textRange, [int]$index, [string]$value) {
    $para = $textRange.Paragraphs($index, 1)
    $para.Text = "-"
    $para = $textRange.Paragraphs($index, 1)
    $para.Text = $value
}

# 1. Drop the trailing "References" slide (slide 5).
$p.Slides.Item(5).Delete()

# 2. Slide 1 - title slide.
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "The Future of AI"
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Generated by AutoPPT"

# 3. Slide 2 - Defining Artificial Intelligence.
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Defining Artificial Intelligence"
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr2 2 "AI is the simulation of human intelligence by machines"
Set-ParaText $tr2 3 "Involves learning, reasoning, and problem-solving"
Set-ParaText $tr2 4 "Transitioning from Narrow AI to General AI (AGI)"

# 4. Slide 3 - Current Trends in AI.
$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Current Trends in AI"
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr3 2 "Explosion of Generative AI and Large Language Models"
Set-ParaText $tr3 3 "AI integration in healthcare, finance, and engineering"
Set-ParaText $tr3 4 "Real-time translation and advanced multimodal capabilities"

# 5. Slide 4 - Societal Impact & Ethics.
$s4 = $p.Slides.Item(4)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Societal Impact & Ethics"
$tr4 = $s4.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr4 2 "Concerns over job displacement and workforce evolution"
Set-ParaText $tr4 3 "Ethical considerations in AI decision-making (Bias/Transparency)"
Set-ParaText $tr4 4 "The importance of AI alignment and safety protocols"

# 6. Speaker notes for the first three slides.
$s1.NotesPage.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "Introduction to the core concept of AI."
$s2.NotesPage.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "Discussing how AI is currently shaping the world."
$s3.NotesPage.Shapes.Placeholders.Item(2).TextFrame.TextRange.Text = "A look at the broader implications of AI adoption."

"done"
